# Fruta / hortaliza, semanal
# The data rows (2-33) of the single sheet got reshuffled: every row's full
# set of values (Fecha, Calidad, Volumen, Precios, etc.) moved to a new row
# position - this is a pure permutation of the existing rows, no cell
# content was actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 33
$firstCol = 1   # A
$lastCol = 20   # T

# Snapshot all the existing row values before writing anything back, so that
# overlapping writes never clobber data we still need to read.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$original = $srcRange.Value2

# Maps the destination row number -> source row number (both in the
# original, pre-edit layout).
$rowMap = @{
    2  = 15
    3  = 16
    4  = 21
    5  = 22
    6  = 23
    7  = 13
    8  = 14
    9  = 26
    10 = 27
    11 = 12
    12 = 7
    13 = 8
    14 = 9
    15 = 28
    16 = 29
    17 = 30
    18 = 33
    19 = 2
    20 = 3
    21 = 4
    22 = 5
    23 = 6
    24 = 19
    25 = 20
    26 = 17
    27 = 18
    28 = 10
    29 = 11
    30 = 31
    31 = 32
    32 = 24
    33 = 25
}

# Build the new block of values in memory (rows ordered 2..33) then write it
# back in one shot. $original (read via Value2) is 1-based, 1..numRows and
# 1..numCols. The array we hand back to Value2 for the write must exactly
# match the destination range size and is 0-based (New-Object default).
$numRows = $lastRow - $firstRow + 1
$numCols = $lastCol - $firstCol + 1
$newValues = New-Object 'object[,]' $numRows, $numCols

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow        # 0-based
    $srcIdx = $srcRow - $firstRow + 1      # 1-based (into $original)
    for ($col = 1; $col -le $numCols; $col++) {
        $newValues[$destIdx, ($col - 1)] = $original[$srcIdx, $col]
    }
}

$destRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$destRange.Value2 = $newValues
